{"js": "// Replace the date in the title and the 25 division problems in the table\n// with the new values, preserving all formatting (fonts, sizes, etc.).\nconst replacements = [\n  [\"2025-05-16 Friday\", \"2025-05-17 Saturday\"],\n  [\"15\u00f73=\", \"10\u00f72=\"],\n  [\"39\u00f72=\", \"28\u00f75=\"],\n  [\"57\u00f76=\", \"42\u00f75=\"],\n  [\"60\u00f74=\", \"64\u00f77=\"],\n  [\"13\u00f77=\", \"68\u00f76=\"],\n  [\"61\u00f76=\", \"66\u00f79=\"],\n  [\"51\u00f75=\", \"62\u00f75=\"],\n  [\"43\u00f79=\", \"73\u00f74=\"],\n  [\"45\u00f79=\", \"60\u00f72=\"],\n  [\"86\u00f76=\", \"40\u00f78=\"],\n  [\"87\u00f75=\", \"19\u00f72=\"],\n  [\"89\u00f73=\", \"64\u00f78=\"],\n  [\"42\u00f72=\", \"64\u00f77=\"],\n  [\"10\u00f75=\", \"69\u00f73=\"],\n  [\"21\u00f75=\", \"97\u00f72=\"],\n  [\"29\u00f76=\", \"56\u00f75=\"],\n  [\"52\u00f72=\", \"41\u00f77=\"],\n  [\"71\u00f73=\", \"38\u00f76=\"],\n  [\"26\u00f76=\", \"71\u00f79=\"],\n  [\"43\u00f78=\", \"16\u00f74=\"],\n  [\"63\u00f76=\", \"40\u00f77=\"],\n  [\"84\u00f79=\", \"83\u00f75=\"],\n  [\"38\u00f72=\", \"58\u00f73=\"],\n  [\"62\u00f73=\", \"70\u00f77=\"],\n  [\"96\u00f73=\", \"86\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date in the title and the 25 division problems in the table\n# with the new values, preserving all formatting (fonts, sizes, etc.).\n$replacements = @(\n  @{old=\"2025-05-16 Friday\"; new=\"2025-05-17 Saturday\"},\n  @{old=\"15\u00f73=\"; new=\"10\u00f72=\"},\n  @{old=\"39\u00f72=\"; new=\"28\u00f75=\"},\n  @{old=\"57\u00f76=\"; new=\"42\u00f75=\"},\n  @{old=\"60\u00f74=\"; new=\"64\u00f77=\"},\n  @{old=\"13\u00f77=\"; new=\"68\u00f76=\"},\n  @{old=\"61\u00f76=\"; new=\"66\u00f79=\"},\n  @{old=\"51\u00f75=\"; new=\"62\u00f75=\"},\n  @{old=\"43\u00f79=\"; new=\"73\u00f74=\"},\n  @{old=\"45\u00f79=\"; new=\"60\u00f72=\"},\n  @{old=\"86\u00f76=\"; new=\"40\u00f78=\"},\n  @{old=\"87\u00f75=\"; new=\"19\u00f72=\"},\n  @{old=\"89\u00f73=\"; new=\"64\u00f78=\"},\n  @{old=\"42\u00f72=\"; new=\"64\u00f77=\"},\n  @{old=\"10\u00f75=\"; new=\"69\u00f73=\"},\n  @{old=\"21\u00f75=\"; new=\"97\u00f72=\"},\n  @{old=\"29\u00f76=\"; new=\"56\u00f75=\"},\n  @{old=\"52\u00f72=\"; new=\"41\u00f77=\"},\n  @{old=\"71\u00f73=\"; new=\"38\u00f76=\"},\n  @{old=\"26\u00f76=\"; new=\"71\u00f79=\"},\n  @{old=\"43\u00f78=\"; new=\"16\u00f74=\"},\n  @{old=\"63\u00f76=\"; new=\"40\u00f77=\"},\n  @{old=\"84\u00f79=\"; new=\"83\u00f75=\"},\n  @{old=\"38\u00f72=\"; new=\"58\u00f73=\"},\n  @{old=\"62\u00f73=\"; new=\"70\u00f77=\"},\n  @{old=\"96\u00f73=\"; new=\"86\u00f73=\"}\n)\n\n$d = $word.ActiveDocument\n\nforeach ($p in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Execute($p.old, $false, $false, $false, $false, $false, $true, 1, $false, $p.new, 2) | Out-Null\n}\n"}
